$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the old sub-header row (row 2: Hiver/Eté/Année/(MW)/(GWh) labels).
# This shifts every data row below it up by one.
$ws.Rows.Item(2).Delete()

# Build the single consolidated header row 1.
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# A1:D1 are brand new headers -> plain/default formatting.
$ws.Range("A1:D1").Style = "Normal"

# E1 carried over the old "(MW)" style (fontId 1) from before the row
# delete; the new layout wants it back to the plain default style.
$ws.Range("E1").Style = "Normal"

# F1:K1 use the font-only header style (Arial 9, General number format).
$headerRange = $ws.Range("F1:K1")
$headerRange.Font.Name = "Arial"
$headerRange.Font.Size = 9

$ws.Range("A2:K2").Select()
